$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used throughout for "percent-like" text values: assigning
# a string such as "84.9%" straight to .Value gets auto-coerced by Excel
# into a number with a brand-new percentage number format (adds a style).
# Writing it as a ="text" formula and then pasting-special as values-only
# collapses it back to a plain text cell using the ORIGINAL style/format.
function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

# --- Column I width: 14 -> 10 ---
$ws.Columns.Item(9).ColumnWidth = 10

# --- Summary box (K/L) updates ---
$ws.Range("L6").Value = 135
$ws.Range("L7").Value = 0
Set-TextValue "L9" "84.9%"
Set-TextValue "L10" "71.2%"

# --- "System, X" -> "X, System" swaps in column G ---
$ws.Range("G11").Value = "dnasr281@gmail.com, System"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("G44").Value = "dnasr281@gmail.com, System"
$ws.Range("G65").Value = "dnasr281@gmail.com, System"
$ws.Range("G71").Value = "dnasr281@gmail.com, System"
$ws.Range("G96").Value = "dnasr281@gmail.com, System"
$ws.Range("G97").Value = "dnasr281@gmail.com, System"
$ws.Range("G99").Value = "dnasr281@gmail.com, System"
$ws.Range("G122").Value = "dnasr281@gmail.com, System"
$ws.Range("G123").Value = "dnasr281@gmail.com, System"
$ws.Range("G125").Value = "dnasr281@gmail.com, System"
$ws.Range("G148").Value = "dnasr281@gmail.com, System"
$ws.Range("G149").Value = "dnasr281@gmail.com, System"
$ws.Range("G151").Value = "dnasr281@gmail.com, System"

# --- Group stats block (rows 18-20) ---
$ws.Range("O18").Value = 21
$ws.Range("P18").Value = 0
Set-TextValue "R18" "80.8%"
Set-TextValue "S18" "76.5%"

$ws.Range("O19").Value = 21
$ws.Range("P19").Value = 0
Set-TextValue "R19" "80.8%"
Set-TextValue "S19" "74.2%"

$ws.Range("O20").Value = 21
$ws.Range("P20").Value = 0
Set-TextValue "R20" "80.8%"
Set-TextValue "S20" "81.3%"

# --- Rows 103, 129, 155: now-recorded sessions ---
# Restyle A:I from s="9" to s="2" (reuse existing style, don't create a new one)
$ws.Range("A2:I2").Copy()
$ws.Range("A103:I103").PasteSpecial(-4122)
$ws.Range("A2:I2").Copy()
$ws.Range("A129:I129").PasteSpecial(-4122)
$ws.Range("A2:I2").Copy()
$ws.Range("A155:I155").PasteSpecial(-4122)

$ws.Range("G103").Value = "dnasr281@gmail.com"
$ws.Range("H103").Value = "36/56"
$ws.Range("I103").Value = "Recorded"

$ws.Range("G129").Value = "dnasr281@gmail.com"
$ws.Range("H129").Value = "28/55"
$ws.Range("I129").Value = "Recorded"

$ws.Range("G155").Value = "dnasr281@gmail.com"
$ws.Range("H155").Value = "32/57"
$ws.Range("I155").Value = "Recorded"
